$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 160.33333
$ws.Range("I9").Value = 196.25
$ws.Range("J9").Value = 88.5
$ws.Range("K9").Value = 196.25
$ws.Range("L9").Value = 88.5
$ws.Range("M9").Value = -27.25
$ws.Range("N9").Value = -426.5
$ws.Range("H41").Value = 491.27274
$ws.Range("I41").Value = 143.71428
$ws.Range("J41").Value = 1099.5
$ws.Range("K41").Value = 143.71428
$ws.Range("L41").Value = 1099.5
$ws.Range("M41").Value = 296.28572
$ws.Range("N41").Value = -1979.5
$ws.Range("H53").Value = 187.83333
$ws.Range("J53").Value = 170.75
$ws.Range("L53").Value = 170.75
$ws.Range("N53").Value = -1444.75
$ws.Range("H99").Value = 535.5
$ws.Range("I99").Value = 535.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1606.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -108.5
$ws.Range("N99").ClearContents()
$ws.Range("H132").Value = 20849.7
$ws.Range("I132").Value = 20849.7
$ws.Range("K132").Value = 62549.10000000001
$ws.Range("M132").Value = -60019.10000000001
$ws.Range("H138").Value = 4504.385
$ws.Range("J138").Value = 3999
$ws.Range("L138").Value = 11997
$ws.Range("N138").Value = -22277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2400
$ws.Range("J45").Value = 5200
$ws.Range("L45").Value = 5200
$ws.Range("N45").Value = -5954
$ws.Range("H74").Value = 694827.5600000001
$ws.Range("I74").Value = 752683.1
$ws.Range("J74").Value = 23703.2
$ws.Range("K74").Value = 752683.1
$ws.Range("L74").Value = 23703.2
$ws.Range("M74").Value = -751809.1
$ws.Range("N74").Value = -25451.2
$ws.Range("H77").Value = 694827.5600000001
$ws.Range("I77").Value = 752683.1
$ws.Range("J77").Value = 23703.2
$ws.Range("K77").Value = 3763415.5
$ws.Range("L77").Value = 118516
$ws.Range("M77").Value = -3759047.5
$ws.Range("N77").Value = -127252
$ws.Range("H110").Value = 1997.6428
$ws.Range("I110").Value = 1999.6666
$ws.Range("K110").Value = 1999.6666
$ws.Range("M110").Value = 45.33339999999998
$ws.Range("H122").Value = 1463
$ws.Range("I122").Value = 1433.0625
$ws.Range("K122").Value = 4299.1875
$ws.Range("M122").Value = -1849.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2033652.8
$ws.Range("I134").Value = 1173.1794
$ws.Range("J134").Value = 41667004
$ws.Range("K134").Value = 3519.5382
$ws.Range("L134").Value = 125001012
$ws.Range("M134").Value = -984.5382
$ws.Range("N134").Value = -125006082

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 126728
$ws.Range("I16").Value = 2001.8334
$ws.Range("K16").Value = 2001.8334
$ws.Range("M16").Value = -1714.8334
$ws.Range("H22").Value = 1861.4375
$ws.Range("I22").Value = 1411
$ws.Range("J22").Value = 2311.875
$ws.Range("K22").Value = 1411
$ws.Range("L22").Value = 2311.875
$ws.Range("M22").Value = -1061
$ws.Range("N22").Value = -3011.875
$ws.Range("H58").Value = 18795878
$ws.Range("I58").Value = 27781998
$ws.Range("K58").Value = 27781998
$ws.Range("M58").Value = -27781795
$ws.Range("H97").Value = 23000
$ws.Range("I97").Value = 23000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 23000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -22009
$ws.Range("N97").ClearContents()
$ws.Range("H99").Value = 22066.4
$ws.Range("I99").Value = 30337.715
$ws.Range("J99").Value = 2766.6667
$ws.Range("K99").Value = 30337.715
$ws.Range("L99").Value = 2766.6667
$ws.Range("M99").Value = -28839.715
$ws.Range("N99").Value = -5762.6667
$ws.Range("H113").Value = 126728
$ws.Range("I113").Value = 2001.8334
$ws.Range("K113").Value = 2001.8334
$ws.Range("M113").Value = 168.1666
$ws.Range("H126").Value = 22066.4
$ws.Range("I126").Value = 30337.715
$ws.Range("J126").Value = 2766.6667
$ws.Range("K126").Value = 91013.145
$ws.Range("L126").Value = 8300.000100000001
$ws.Range("M126").Value = -88543.145
$ws.Range("N126").Value = -13240.0001
$ws.Range("H132").Value = 3593
$ws.Range("I132").Value = 3363.1667
$ws.Range("K132").Value = 10089.5001
$ws.Range("M132").Value = -7559.500100000001
$ws.Range("H136").Value = 18795878
$ws.Range("I136").Value = 27781998
$ws.Range("K136").Value = 83345994
$ws.Range("M136").Value = -83343444

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 210.36363
$ws.Range("I26").Value = 151.4
$ws.Range("J26").Value = 800
$ws.Range("K26").Value = 454.2
$ws.Range("L26").Value = 2400
$ws.Range("M26").Value = -166.2
$ws.Range("N26").Value = -2976
$ws.Range("H33").Value = 3000.3333
$ws.Range("I33").Value = 1
$ws.Range("K33").Value = 6
$ws.Range("M33").Value = 277
$ws.Range("H60").Value = 3022.4285
$ws.Range("I60").Value = 3022.4285
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 9067.2855
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -8816.2855
$ws.Range("N60").ClearContents()
$ws.Range("H87").Value = 13931
$ws.Range("I87").Value = 5370.625
$ws.Range("J87").Value = 23714.285
$ws.Range("K87").Value = 16111.875
$ws.Range("L87").Value = 71142.855
$ws.Range("M87").Value = -14863.875
$ws.Range("N87").Value = -73638.855
$ws.Range("H90").Value = 13931
$ws.Range("I90").Value = 5370.625
$ws.Range("J90").Value = 23714.285
$ws.Range("K90").Value = 48335.625
$ws.Range("L90").Value = 213428.565
$ws.Range("M90").Value = -42095.625
$ws.Range("N90").Value = -225908.565
$ws.Range("H121").Value = 3233.55
$ws.Range("I121").Value = 1176.6666
$ws.Range("J121").Value = 3596.5293
$ws.Range("K121").Value = 3529.9998
$ws.Range("L121").Value = 10789.5879
$ws.Range("M121").Value = -2219.9998
$ws.Range("N121").Value = -13409.5879
$ws.Range("H129").Value = 627573.2
$ws.Range("I129").Value = 1251802.4
$ws.Range("J129").Value = 3344
$ws.Range("K129").Value = 3755407.2
$ws.Range("L129").Value = 10032
$ws.Range("M129").Value = -3750407.2
$ws.Range("N129").Value = -20032
$ws.Range("H131").Value = 4398.347
$ws.Range("I131").Value = 546.1429000000001
$ws.Range("J131").Value = 5939.2285
$ws.Range("K131").Value = 1638.4287
$ws.Range("L131").Value = 17817.6855
$ws.Range("M131").Value = 3401.5713
$ws.Range("N131").Value = -27897.6855
$ws.Range("H140").Value = 2844.923
$ws.Range("I140").Value = 2332
$ws.Range("J140").Value = 9000
$ws.Range("K140").Value = 6996
$ws.Range("L140").Value = 27000
$ws.Range("M140").Value = -1816
$ws.Range("N140").Value = -37360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9453.657999999999
$ws.Range("I70").Value = 10835.9
$ws.Range("J70").Value = 8137.2383
$ws.Range("K70").Value = 10835.9
$ws.Range("L70").Value = 8137.2383
$ws.Range("M70").Value = -10565.9
$ws.Range("N70").Value = -8677.238300000001
$ws.Range("H73").Value = 9453.657999999999
$ws.Range("I73").Value = 10835.9
$ws.Range("J73").Value = 8137.2383
$ws.Range("K73").Value = 10835.9
$ws.Range("L73").Value = 8137.2383
$ws.Range("M73").Value = -9899.9
$ws.Range("N73").Value = -10009.2383
$ws.Range("H102").Value = 1691.5454
$ws.Range("I102").Value = 1623
$ws.Range("K102").Value = 1623
$ws.Range("M102").Value = -1
$ws.Range("H113").Value = 1370.375
$ws.Range("I113").Value = 1370.375
$ws.Range("K113").Value = 1370.375
$ws.Range("M113").Value = 799.625
$ws.Range("H122").Value = 2996.12
$ws.Range("I122").Value = 2961.1304
$ws.Range("J122").Value = 3398.5
$ws.Range("K122").Value = 8883.3912
$ws.Range("L122").Value = 10195.5
$ws.Range("M122").Value = -6433.3912
$ws.Range("N122").Value = -15095.5
$ws.Range("H126").Value = 2589.5
$ws.Range("I126").Value = 2134.5
$ws.Range("K126").Value = 6403.5
$ws.Range("M126").Value = -3933.5
$ws.Range("H132").Value = 8919.937
$ws.Range("I132").Value = 5164.8
$ws.Range("K132").Value = 15494.4
$ws.Range("M132").Value = -12964.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3524.875
$ws.Range("I7").Value = 2885.5715
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 2885.5715
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -2773.5715
$ws.Range("N7").Value = -8224
$ws.Range("H40").Value = 4032.8125
$ws.Range("J40").Value = 5194.6
$ws.Range("L40").Value = 5194.6
$ws.Range("N40").Value = -5466.6
$ws.Range("H46").Value = 4260.316
$ws.Range("J46").Value = 5614.769
$ws.Range("L46").Value = 5614.769
$ws.Range("N46").Value = -5990.769
$ws.Range("H55").Value = 1429.0883
$ws.Range("I55").Value = 1276.3125
$ws.Range("K55").Value = 1276.3125
$ws.Range("M55").Value = -1103.3125
$ws.Range("H122").Value = 3523
$ws.Range("J122").Value = 4330.1665
$ws.Range("L122").Value = 12990.4995
$ws.Range("N122").Value = -17890.4995
$ws.Range("H126").Value = 3524.875
$ws.Range("I126").Value = 2885.5715
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 8656.7145
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -6186.7145
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 1013233.1
$ws.Range("I132").Value = 1669547.1
$ws.Range("K132").Value = 5008641.300000001
$ws.Range("M132").Value = -5006111.300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 503501.5
$ws.Range("I4").Value = 1000000
$ws.Range("J4").Value = 7003
$ws.Range("K4").Value = 1000000
$ws.Range("L4").Value = 7003
$ws.Range("M4").Value = -999887
$ws.Range("N4").Value = -7229
$ws.Range("H107").Value = 1024.6428
$ws.Range("J107").Value = 1333.3334
$ws.Range("L107").Value = 4000.0002
$ws.Range("N107").Value = -7840.0002
$ws.Range("H122").Value = 77679.47
$ws.Range("I122").Value = 2988.1
$ws.Range("K122").Value = 8964.299999999999
$ws.Range("M122").Value = -6514.299999999999
$ws.Range("H126").Value = 1600
$ws.Range("I126").Value = 1250
$ws.Range("J126").Value = 1833.3334
$ws.Range("K126").Value = 3750
$ws.Range("L126").Value = 5500.0002
$ws.Range("M126").Value = -1280
$ws.Range("N126").Value = -10440.0002
